$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "03-" + "31" -> single run "03-31-2023 ---...".
#    The paragraph currently reads "03-" | "31" | "-2023 ------...." across
#    three runs; collapse the first two into one run so the date renders as
#    a single "03-31-2023" token, matching the committed text exactly.
# ---------------------------------------------------------------------------
$dateFind = $d.Content
$ok1 = $dateFind.Find.Execute("03-31", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "03-31", 2)

# ---------------------------------------------------------------------------
# 2) Append the new log entries (three new paragraphs) right after the last
#    paragraph in the document ("Results somewhat similar ... episodes).").
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.InsertParagraphAfter()

$p6 = $d.Paragraphs($d.Paragraphs.Count)
$p6.Range.InsertAfter("Next attempt reduce runtime closer to previous (40s). Increase x and y bounds to 35.")
$p6.Range.InsertParagraphAfter()

$p7 = $d.Paragraphs($d.Paragraphs.Count)
$p7.Range.InsertAfter("04-04-2023 ---------------------------------------------------")
$p7.Range.InsertParagraphAfter()

$p8 = $d.Paragraphs($d.Paragraphs.Count)
$p8.Range.InsertAfter("Similar results. Not much improvement. Possibly maximizing the system given bounds ie needs more room roll over and come back up.")
